# Update the "Förändrad" date column (C) for all data rows (2-61)
# from 45189 (2023-09-20) to 45190 (2023-09-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 61; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
